# Auto update Excel log
# Appends new sensor-event rows to the Proximity, mmWave and Camera sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Proximity sheet: add rows 10-11 (Living Room Main Door ENTER/EXIT)
# ---------------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "15:53:12", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:53:50", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$r = 10
foreach ($row in $proximityRows) {
    $wsProximity.Cells.Item($r, 1).Value = $row[0]
    $wsProximity.Cells.Item($r, 2).Value = $row[1]
    $wsProximity.Cells.Item($r, 3).Value = $row[2]
    $wsProximity.Cells.Item($r, 4).Value = $row[3]
    $wsProximity.Cells.Item($r, 5).Value = $row[4]
    $wsProximity.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# mmWave sheet: add rows 11-16 (Living Room PRESENCE_DETECTED / Active)
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$mmWaveTimestamps = @("15:53:09", "15:53:19", "15:53:30", "15:53:40", "15:53:51", "15:54:01")

$r = 11
foreach ($ts in $mmWaveTimestamps) {
    $wsMmWave.Cells.Item($r, 1).Value = "2026-02-01"
    $wsMmWave.Cells.Item($r, 2).Value = $ts
    $wsMmWave.Cells.Item($r, 3).Value = "15:00"
    $wsMmWave.Cells.Item($r, 4).Value = "Living Room"
    $wsMmWave.Cells.Item($r, 5).Value = "PRESENCE_DETECTED"
    $wsMmWave.Cells.Item($r, 6).Value = "Active"
    $r++
}

# ---------------------------------------------------------------------------
# Camera sheet: add rows 8-9 (Living Room Main Door Image Captured / Active)
# ---------------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")

$cameraTimestamps = @("15:53:13", "15:53:50")

$r = 8
foreach ($ts in $cameraTimestamps) {
    $wsCamera.Cells.Item($r, 1).Value = "2026-02-01"
    $wsCamera.Cells.Item($r, 2).Value = $ts
    $wsCamera.Cells.Item($r, 3).Value = "15:00"
    $wsCamera.Cells.Item($r, 4).Value = "Living Room Main Door"
    $wsCamera.Cells.Item($r, 5).Value = "Image Captured"
    $wsCamera.Cells.Item($r, 6).Value = "Active"
    $r++
}
